# "reset date button init"
#
# Swap the two ticker rows' date/label values and update the price figures.
# Row 1 (was 2023-3-12 / A / 127.58) -> 2023-3-1 / D / 79.29
# Row 2 (was 2023-3-12 / D / 81.32)  -> 2023-3-1 / A / 131.93
#
# These cells are plain text (not real dates/numbers), but values like
# "2023-3-1" and "79.29" look numeric/date-like to Excel's auto-detection,
# so a leading apostrophe forces them to stay text (the apostrophe itself
# is not part of the stored value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "'2023-3-1"
$ws.Range("B1").Value = "'D"
$ws.Range("D1").Value = "'79.29`n"

$ws.Range("A2").Value = "'2023-3-1"
$ws.Range("B2").Value = "'A"
$ws.Range("D2").Value = "'131.93`n"
